$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.248.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "'2.578.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'555.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "'142.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "'2.584.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "'0.165"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.07%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'3.032.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "'59.225.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "'23.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'2.585.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'337.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'10.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "'6.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.473"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.87%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'62.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.54%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").Value = "'7.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'0.0₃0774"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'6.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "'158.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'19.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "'4.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "'0.896"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "'37.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.853"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "'289.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").Value = "'137.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.84%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").Value = "'0.593"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "'0.0530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "'18.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "'1.938.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.13%  "
